$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column Q (day 16) values for rows 2-6
$ws.Range("Q2").Value = 18261.4
$ws.Range("Q3").Value = 4882
$ws.Range("Q4").Value = 2772
$ws.Range("Q5").Value = 2599
$ws.Range("Q6").Value = 28514.4

# Update column AG (total) values for rows 2-6
$ws.Range("AG2").Value = 185911.26
$ws.Range("AG3").Value = 80407.8
$ws.Range("AG4").Value = 50341.69
$ws.Range("AG5").Value = 49004.2
$ws.Range("AG6").Value = 365664.95
